$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = 0.2560897945784967
$ws.Range("C5").Value = 0.1922302243704833
$ws.Range("D5").Value = 0.1453653151772641
$ws.Range("E5").Value = 0.1929223575287684
$ws.Range("F5").Value = 0.2240464861714938
$ws.Range("G5").Value = 0.2544840169573067

$ws.Range("B6").Value = 0.2770267008588849
$ws.Range("C6").Value = 0.1605440980871795
$ws.Range("D6").Value = 0.1369792959557058
$ws.Range("E6").Value = 0.208388836430242
$ws.Range("F6").Value = 0.2068473773143613
$ws.Range("G6").Value = 0.2132499575263846

$ws.Range("B7").Value = 0.1748995270691545
$ws.Range("C7").Value = 0.08298676030387352
$ws.Range("D7").Value = 0.09247763512786311
$ws.Range("E7").Value = 0.1442451350674593
$ws.Range("F7").Value = 0.1477468036115866
$ws.Range("G7").Value = 0.1175250756782335

$ws.Range("B8").Value = 0.1396381204092653
$ws.Range("C8").Value = 0.03970428935642073
$ws.Range("D8").Value = 0.05091816180503488
$ws.Range("E8").Value = 0.1933786134269216
$ws.Range("F8").Value = 0.1249432911345317
$ws.Range("G8").Value = 0.1813139052758002
